$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove existing hyperlinks first; Insert() does not reliably relocate
# hyperlinks that live below the inserted block, so we recreate every
# hyperlink afterwards at its correct final address instead.
$ws.Cells.Hyperlinks.Delete()

# Insert 8 new rows before row 15 (shifts old rows 16-18 down to 24-26)
$ws.Rows("15:22").Insert()

# New credit rows (freesound.org sound effects) added to the bottom of the
# existing table section
$newRows = @(
    @{ Row = 15; Name = "60013__qubodup__whoosh"; Url = "https://freesound.org/people/qubodup/sounds/60013/" },
    @{ Row = 16; Name = "382735__schots__gun-shot"; Url = "https://freesound.org/people/schots/sounds/382735/" },
    @{ Row = 17; Name = "588246__rkkaleikau__energy-weapon-laser"; Url = "https://freesound.org/people/rkkaleikau/sounds/588246/" },
    @{ Row = 18; Name = "566435__merrick079__punch2"; Url = "https://freesound.org/people/Merrick079/sounds/566435/" },
    @{ Row = 19; Name = "232358__richerlandtv__heavy-impacts"; Url = "https://freesound.org/people/RICHERlandTV/sounds/232358/" },
    @{ Row = 20; Name = "341247__sharesynth__jump01"; Url = "https://freesound.org/people/sharesynth/sounds/341247/" },
    @{ Row = 21; Name = "561646__mattruthsound__hit-punch-cloth-pillow-bedding-004"; Url = "https://freesound.org/people/MattRuthSound/sounds/561646/" },
    @{ Row = 22; Name = "433644__dersuperanton__game-over-sound"; Url = "https://freesound.org/people/dersuperanton/sounds/433644/" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value2 = $r.Name
    $ws.Range("B$rowNum").Value2 = $r.Url
    $ws.Range("C$rowNum").Value2 = "Doesn't require attribution"
    $ws.Range("C$rowNum").Style = "Normal"
}

# Re-create every hyperlink in the sheet at its final resting cell
$allLinks = @(
    @{ Cell = "B2";  Url = "https://freesound.org/people/Whiprealgood/sounds/87535/" },
    @{ Cell = "B3";  Url = "https://freesound.org/people/suntemple/sounds/253172/" },
    @{ Cell = "B4";  Url = "https://www.fontspace.com/a-area-kilometer-50-font-f53888" },
    @{ Cell = "B5";  Url = "https://opengameart.org/content/simple-explosion-bleeds-game-art" },
    @{ Cell = "B6";  Url = "https://opengameart.org/content/spikes-0" },
    @{ Cell = "B7";  Url = "https://opengameart.org/content/various-inventory-24-pixel-icon-set" },
    @{ Cell = "B8";  Url = "https://opengameart.org/content/energy-icon" },
    @{ Cell = "B15"; Url = "https://freesound.org/people/qubodup/sounds/60013/" },
    @{ Cell = "B16"; Url = "https://freesound.org/people/schots/sounds/382735/" },
    @{ Cell = "B17"; Url = "https://freesound.org/people/rkkaleikau/sounds/588246/" },
    @{ Cell = "B18"; Url = "https://freesound.org/people/Merrick079/sounds/566435/" },
    @{ Cell = "B19"; Url = "https://freesound.org/people/RICHERlandTV/sounds/232358/" },
    @{ Cell = "B20"; Url = "https://freesound.org/people/sharesynth/sounds/341247/" },
    @{ Cell = "B21"; Url = "https://freesound.org/people/MattRuthSound/sounds/561646/" },
    @{ Cell = "B22"; Url = "https://freesound.org/people/dersuperanton/sounds/433644/" },
    @{ Cell = "B25"; Url = "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack" },
    @{ Cell = "B26"; Url = "https://free-game-assets.itch.io/night-city-street-2d-background-tiles" }
)

foreach ($l in $allLinks) {
    $ws.Hyperlinks.Add($ws.Range($l.Cell), $l.Url) | Out-Null
    $ws.Range($l.Cell).Style = "Hyperlink"
}

# Restore the scroll position / selection recorded by Excel on save
$wb.Windows.Item(1).ScrollRow = 16
$ws.Range("C30").Select()
